$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "All Sites": in-place value corrections + a few brand-new cells in
# existing rows (no row shifting here).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Sites")

# Row 12 - ID 51 was a typo for 11; species record flips Present -> Absent
$wsAll.Cells.Item(12, 1).Value = 11
$wsAll.Cells.Item(12, 2).Value = "Absent"

# Row 22 - corrected counts
$wsAll.Cells.Item(22, 8).Value = 4
$wsAll.Cells.Item(22, 11).Value = 1
$wsAll.Cells.Item(22, 12).Value = 13

# Row 23 - corrected counts
$wsAll.Cells.Item(23, 6).Value = 8
$wsAll.Cells.Item(23, 11).Value = 2
$wsAll.Cells.Item(23, 14).Value = 6

# Row 54 - Absent -> Present
$wsAll.Cells.Item(54, 2).Value = "Present"

# Row 92 - corrected counts
$wsAll.Cells.Item(92, 5).Value = 2
$wsAll.Cells.Item(92, 13).Value = 3
$wsAll.Cells.Item(92, 16).Value = 5

# Rows 113-116 - newly recorded counts for Efr (G) / Dm (L)
$wsAll.Cells.Item(113, 7).Value = 4
$wsAll.Cells.Item(114, 7).Value = 10
$wsAll.Cells.Item(114, 12).Value = 1
$wsAll.Cells.Item(115, 7).Value = 1
$wsAll.Cells.Item(116, 12).Value = 1

# ---------------------------------------------------------------------------
# Sheet "Tl": insert one new site row (92) before the existing row 12,
# pushing the rest down.
# ---------------------------------------------------------------------------
$wsTl = $wb.Worksheets.Item("Tl")
$wsTl.Rows(12).Insert()
$wsTl.Cells.Item(12, 1).Value = 92
$wsTl.Cells.Item(12, 2).Value = 30.02547222222222
$wsTl.Cells.Item(12, 3).Value = -90.115638888888881
$wsTl.Cells.Item(12, 4).Value = "Present"
# A/B/C on the source row above carried no explicit style - match that.
$wsTl.Range("A12:C12").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "Th": insert one new site row (111) before existing row 28, then
# append two more new rows (112, 113) at the end.
# ---------------------------------------------------------------------------
$wsTh = $wb.Worksheets.Item("Th")
$wsTh.Rows(28).Insert()
$wsTh.Cells.Item(28, 1).Value = 111
$wsTh.Cells.Item(28, 2).Value = 29.785
$wsTh.Cells.Item(28, 3).Value = -90.407333333333341
$wsTh.Cells.Item(28, 4).Value = "Present"
$wsTh.Range("A28:C28").Style = "Normal"

$wsTh.Cells.Item(30, 1).Value = 112
$wsTh.Cells.Item(30, 2).Value = 29.824071666666665
$wsTh.Cells.Item(30, 3).Value = -90.476001111111117
$wsTh.Cells.Item(30, 4).Value = "Present"

$wsTh.Cells.Item(31, 1).Value = 113
$wsTh.Cells.Item(31, 2).Value = 29.646777777777778
$wsTh.Cells.Item(31, 3).Value = -90.540944444444449
$wsTh.Cells.Item(31, 4).Value = "Present"

# ---------------------------------------------------------------------------
# Sheet "Sl": insert one new site row (21) before existing row 3.
# ---------------------------------------------------------------------------
$wsSl = $wb.Worksheets.Item("Sl")
$wsSl.Rows(3).Insert()
$wsSl.Cells.Item(3, 1).Value = 21
$wsSl.Cells.Item(3, 2).Value = 32.32
$wsSl.Cells.Item(3, 3).Value = -93.67
$wsSl.Cells.Item(3, 4).Value = "Present"

# ---------------------------------------------------------------------------
# Sheet "Efr": append two new rows at the bottom.
# ---------------------------------------------------------------------------
$wsEfr = $wb.Worksheets.Item("Efr")
$wsEfr.Cells.Item(32, 1).Value = 112
$wsEfr.Cells.Item(32, 2).Value = 29.824071666666665
$wsEfr.Cells.Item(32, 3).Value = -90.476001111111117
$wsEfr.Cells.Item(32, 4).Value = "Present"
$wsEfr.Cells.Item(32, 4).HorizontalAlignment = -4152

$wsEfr.Cells.Item(33, 1).Value = 114
$wsEfr.Cells.Item(33, 2).Value = 29.86675
$wsEfr.Cells.Item(33, 3).Value = -90.599666666666664
$wsEfr.Cells.Item(33, 4).Value = "Present"
$wsEfr.Cells.Item(33, 4).HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# Sheet "Efl": append one new row at the bottom.
# ---------------------------------------------------------------------------
$wsEfl = $wb.Worksheets.Item("Efl")
$wsEfl.Cells.Item(4, 1).Value = 92
$wsEfl.Cells.Item(4, 2).Value = 30.02547222222222
$wsEfl.Cells.Item(4, 3).Value = -90.115638888888881
$wsEfl.Cells.Item(4, 4).Value = "Present"

# ---------------------------------------------------------------------------
# View state: selections / active sheet, mirroring the saved workbook state.
# ---------------------------------------------------------------------------
[void]$wsTl.Range("A1").Select()

[void]$wsAll.Range("A115:D115").Select()

[void]$wsSl.Range("G12").Select()

[void]$wsEfr.Range("F33").Select()

[void]$wsEfl.Range("C6").Select()

$wsTh.Activate()
[void]$wsTh.Range("D31").Select()

Write-Output "edits applied"
